$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3372.157
$ws.Range("I64").Value = 3109.2683
$ws.Range("J64").Value = 4450
$ws.Range("K64").Value = 3109.2683
$ws.Range("L64").Value = 4450
$ws.Range("M64").Value = -2861.2683
$ws.Range("N64").Value = -4946

$ws.Range("H67").Value = 3372.157
$ws.Range("I67").Value = 3109.2683
$ws.Range("J67").Value = 4450
$ws.Range("K67").Value = 3109.2683
$ws.Range("L67").Value = 4450
$ws.Range("M67").Value = -2251.2683
$ws.Range("N67").Value = -6166

$ws.Range("H129").Value = 853.81396
$ws.Range("I129").Value = 512.125
$ws.Range("J129").Value = 888.85895
$ws.Range("K129").Value = 1536.375
$ws.Range("L129").Value = 2666.57685
$ws.Range("M129").Value = 3463.625
$ws.Range("N129").Value = -12666.57685

$ws.Range("H137").Value = 1279
$ws.Range("I137").Value = 1375.9333
$ws.Range("J137").Value = 1071.2858
$ws.Range("K137").Value = 4127.7999
$ws.Range("L137").Value = 3213.8574
$ws.Range("M137").Value = -1577.7999
$ws.Range("N137").Value = -8313.857400000001

$ws.Range("H138").Value = 3947.7979
$ws.Range("I138").Value = 2242.6943
$ws.Range("J138").Value = 5006.1377
$ws.Range("K138").Value = 6728.0829
$ws.Range("L138").Value = 15018.4131
$ws.Range("M138").Value = -1588.0829
$ws.Range("N138").Value = -25298.4131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1517.4872
$ws.Range("I74").Value = 1071.4615
$ws.Range("J74").Value = 2409.5386
$ws.Range("K74").Value = 1071.4615
$ws.Range("L74").Value = 2409.5386
$ws.Range("M74").Value = -197.4614999999999
$ws.Range("N74").Value = -4157.5386

$ws.Range("H77").Value = 1517.4872
$ws.Range("I77").Value = 1071.4615
$ws.Range("J77").Value = 2409.5386
$ws.Range("K77").Value = 5357.307499999999
$ws.Range("L77").Value = 12047.693
$ws.Range("M77").Value = -989.307499999999
$ws.Range("N77").Value = -20783.693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3696.7856
$ws.Range("I80").Value = 2588.2727
$ws.Range("J80").Value = 4414.0586
$ws.Range("K80").Value = 2588.2727
$ws.Range("L80").Value = 4414.0586
$ws.Range("M80").Value = -1590.2727
$ws.Range("N80").Value = -6410.0586

$ws.Range("H83").Value = 3696.7856
$ws.Range("I83").Value = 2588.2727
$ws.Range("J83").Value = 4414.0586
$ws.Range("K83").Value = 12941.3635
$ws.Range("L83").Value = 22070.293
$ws.Range("M83").Value = -7949.363499999999
$ws.Range("N83").Value = -32054.293

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3196.23
$ws.Range("I31").Value = 1944.1428
$ws.Range("J31").Value = 4399.216
$ws.Range("K31").Value = 1944.1428
$ws.Range("L31").Value = 4399.216
$ws.Range("M31").Value = -1649.1428
$ws.Range("N31").Value = -4989.216

$ws.Range("H34").Value = 3196.23
$ws.Range("I34").Value = 1944.1428
$ws.Range("J34").Value = 4399.216
$ws.Range("K34").Value = 1944.1428
$ws.Range("L34").Value = 4399.216
$ws.Range("M34").Value = -1742.1428
$ws.Range("N34").Value = -4803.216

$ws.Range("H58").Value = 2825.0322
$ws.Range("I58").Value = 2304.8333
$ws.Range("J58").Value = 3153.5789
$ws.Range("K58").Value = 2304.8333
$ws.Range("L58").Value = 3153.5789
$ws.Range("M58").Value = -2101.8333
$ws.Range("N58").Value = -3559.5789

$ws.Range("H86").Value = 38467640
$ws.Range("I86").Value = 55560756
$ws.Range("J86").Value = 8125
$ws.Range("K86").Value = 55560756
$ws.Range("L86").Value = 8125
$ws.Range("M86").Value = -55559633
$ws.Range("N86").Value = -10371

$ws.Range("H89").Value = 38467640
$ws.Range("I89").Value = 55560756
$ws.Range("J89").Value = 8125
$ws.Range("K89").Value = 277803780
$ws.Range("L89").Value = 40625
$ws.Range("M89").Value = -277798164
$ws.Range("N89").Value = -51857

$ws.Range("H94").Value = 1749.6316
$ws.Range("I94").Value = 2140.4
$ws.Range("J94").Value = 1315.4445
$ws.Range("K94").Value = 2140.4
$ws.Range("L94").Value = 1315.4445
$ws.Range("M94").Value = -1689.4
$ws.Range("N94").Value = -2217.4445

$ws.Range("H99").Value = 1785437.1
$ws.Range("I99").Value = 2463351.5
$ws.Range("J99").Value = 22860
$ws.Range("K99").Value = 2463351.5
$ws.Range("L99").Value = 22860
$ws.Range("M99").Value = -2461853.5
$ws.Range("N99").Value = -25856

$ws.Range("H107").Value = 335.46155
$ws.Range("I107").Value = 282.81818
$ws.Range("K107").Value = 282.81818
$ws.Range("M107").Value = 1637.18182

$ws.Range("H122").Value = 87627.086
$ws.Range("I122").Value = 127063.875
$ws.Range("J122").Value = 8753.5
$ws.Range("K122").Value = 381191.625
$ws.Range("L122").Value = 26260.5
$ws.Range("M122").Value = -378741.625
$ws.Range("N122").Value = -31160.5

$ws.Range("H126").Value = 1785437.1
$ws.Range("I126").Value = 2463351.5
$ws.Range("J126").Value = 22860
$ws.Range("K126").Value = 7390054.5
$ws.Range("L126").Value = 68580
$ws.Range("M126").Value = -7387584.5
$ws.Range("N126").Value = -73520

$ws.Range("H134").Value = 2045.5555
$ws.Range("I134").Value = 2067.0303
$ws.Range("J134").Value = 1809.3334
$ws.Range("K134").Value = 6201.090899999999
$ws.Range("L134").Value = 5428.0002
$ws.Range("M134").Value = -3666.090899999999
$ws.Range("N134").Value = -10498.0002

$ws.Range("H136").Value = 2825.0322
$ws.Range("I136").Value = 2304.8333
$ws.Range("J136").Value = 3153.5789
$ws.Range("K136").Value = 6914.499899999999
$ws.Range("L136").Value = 9460.736699999999
$ws.Range("M136").Value = -4364.499899999999
$ws.Range("N136").Value = -14560.7367

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 744.6429000000001
$ws.Range("I113").Value = 713.7143
$ws.Range("J113").Value = 775.5714
$ws.Range("K113").Value = 2141.1429
$ws.Range("L113").Value = 2326.7142
$ws.Range("M113").Value = 28.85710000000017
$ws.Range("N113").Value = -6666.7142

$ws.Range("H131").Value = 787.65
$ws.Range("I131").Value = 349.0909
$ws.Range("J131").Value = 841.85394
$ws.Range("K131").Value = 1047.2727
$ws.Range("L131").Value = 2525.56182
$ws.Range("M131").Value = 3992.7273
$ws.Range("N131").Value = -12605.56182

$ws.Range("H132").Value = 1858.64
$ws.Range("I132").Value = 884
$ws.Range("J132").Value = 2102.3
$ws.Range("K132").Value = 7956
$ws.Range("L132").Value = 18920.7
$ws.Range("M132").Value = -5426
$ws.Range("N132").Value = -23980.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 50000
$ws.Range("J18").Value = 50000
$ws.Range("L18").Value = 50000
$ws.Range("N18").Value = -50586

$ws.Range("H43").Value = 2886.8
$ws.Range("I43").Value = 2886.8
$ws.Range("K43").Value = 2886.8
$ws.Range("M43").Value = -2735.8

$ws.Range("H57").Value = 18218.334

$ws.Range("H80").Value = 2111.25
$ws.Range("I80").Value = 2023.75
$ws.Range("J80").Value = 2155
$ws.Range("K80").Value = 2023.75
$ws.Range("L80").Value = 2155
$ws.Range("M80").Value = -1025.75
$ws.Range("N80").Value = -4151

$ws.Range("H83").Value = 2111.25
$ws.Range("I83").Value = 2023.75
$ws.Range("J83").Value = 2155
$ws.Range("K83").Value = 10118.75
$ws.Range("L83").Value = 10775
$ws.Range("M83").Value = -5126.75
$ws.Range("N83").Value = -20759

$ws.Range("H107").Value = 1833.8572
$ws.Range("I107").Value = 3576.6667
$ws.Range("J107").Value = 526.75
$ws.Range("K107").Value = 3576.6667
$ws.Range("L107").Value = 526.75
$ws.Range("M107").Value = -1656.6667
$ws.Range("N107").Value = -4366.75

$ws.Range("H113").Value = 2168.2856
$ws.Range("I113").Value = 1872.5
$ws.Range("J113").Value = 2562.6667
$ws.Range("K113").Value = 1872.5
$ws.Range("L113").Value = 2562.6667
$ws.Range("M113").Value = 297.5
$ws.Range("N113").Value = -6902.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1482.909
$ws.Range("I96").Value = 1370.25
$ws.Range("K96").Value = 1370.25
$ws.Range("M96").Value = 2.75

$ws.Range("H107").Value = 471.52173
$ws.Range("I107").Value = 427.06668
$ws.Range("J107").Value = 554.875
$ws.Range("K107").Value = 1281.20004
$ws.Range("L107").Value = 1664.625
$ws.Range("M107").Value = 638.7999599999998
$ws.Range("N107").Value = -5504.625
